# Upload last hw grades (hw13, column N) for students who had recently
# submitted, and fill in a missed hw12 (column M) grade for one student.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# hw13 (column N) grades newly entered
$ws.Range("N3").Value  = 100
$ws.Range("N4").Value  = 95
$ws.Range("N5").Value  = 89
$ws.Range("N7").Value  = 100
$ws.Range("N8").Value  = 100
$ws.Range("N9").Value  = 100
$ws.Range("N10").Value = 100
$ws.Range("N11").Value = 100
$ws.Range("N12").Value = 100
$ws.Range("N13").Value = 90
$ws.Range("N14").Value = 95

# row 15 was missing both hw12 (M) and hw13 (N) grades
$ws.Range("M15").Value = 78
$ws.Range("N15").Value = 90

$ws.Range("N16").Value = 100

# Move the active selection to reflect where editing left off
$ws.Range("M4").Select()

$wb.Save()
